# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-09-08 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-09-09 Tuesday", 2)

# Update the division-problem table. Cells are addressed by (row, column)
# position rather than by old text, because some new values collide with
# other cells' old values elsewhere in the table (e.g. "65÷3=" is both an
# old value in row 17 and the new value produced in row 1), so a global
# text Find/Replace could clobber an already-updated cell.
$t = $d.Tables.Item(1)

function Set-Cell($table, $r, $c, $text) {
    $table.Cell($r, $c).Range.Text = $text
}

# Row 1
Set-Cell $t 1 1 "65÷3="
Set-Cell $t 1 2 "52÷4="
Set-Cell $t 1 3 "50÷7="
Set-Cell $t 1 4 "98÷4="
Set-Cell $t 1 5 "58÷2="

# Row 5
Set-Cell $t 5 1 "61÷2="
Set-Cell $t 5 2 "32÷5="
Set-Cell $t 5 3 "81÷3="
Set-Cell $t 5 4 "32÷2="
Set-Cell $t 5 5 "93÷3="

# Row 9
Set-Cell $t 9 1 "31÷2="
Set-Cell $t 9 2 "72÷4="
Set-Cell $t 9 3 "41÷3="
Set-Cell $t 9 4 "19÷4="
Set-Cell $t 9 5 "18÷9="

# Row 13
Set-Cell $t 13 1 "24÷4="
Set-Cell $t 13 2 "35÷6="
Set-Cell $t 13 3 "66÷8="
Set-Cell $t 13 4 "73÷5="
Set-Cell $t 13 5 "39÷3="

# Row 17
Set-Cell $t 17 1 "48÷7="
Set-Cell $t 17 2 "33÷3="
Set-Cell $t 17 3 "37÷9="
Set-Cell $t 17 4 "54÷2="
Set-Cell $t 17 5 "57÷3="
